$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32 (Leve Item ID 5484)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 2096.4  # H32: was 2283.7778
$ws.Cells.Item(32, 9).Value = 1262.25  # I32: was 1383
$ws.Cells.Item(32, 10).Value = 2652.5  # J32: was 2734.1667
$ws.Cells.Item(32, 11).Value = 1262.25  # K32: was 1383
$ws.Cells.Item(32, 12).Value = 2652.5  # L32: was 2734.1667
$ws.Cells.Item(32, 13).Value = -936.25  # M32: was -1057
$ws.Cells.Item(32, 14).Value = -3304.5  # N32: was -3386.1667

# Sheet ALC, row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 17859082  # H33: was 17859130
$ws.Cells.Item(33, 10).Value = 6076  # J33: was 6300
$ws.Cells.Item(33, 12).Value = 6076  # L33: was 6300
$ws.Cells.Item(33, 14).Value = -6534  # N33: was -6758

# Sheet ALC, row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 14924.75  # H76: was 12899.6
$ws.Cells.Item(76, 10).Value = 25000  # J76: was 18266.334
$ws.Cells.Item(76, 12).Value = 25000  # L76: was 18266.334
$ws.Cells.Item(76, 14).Value = -25630  # N76: was -18896.334

# Sheet ALC, row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 14924.75  # H79: was 12899.6
$ws.Cells.Item(79, 10).Value = 25000  # J79: was 18266.334
$ws.Cells.Item(79, 12).Value = 25000  # L79: was 18266.334
$ws.Cells.Item(79, 14).Value = -27184  # N79: was -20450.334

# Sheet ALC, row 86 (Leve Item ID 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 11399.8  # H86: was 10483.167
$ws.Cells.Item(86, 10).Value = 15250  # J86: was 12133.333
$ws.Cells.Item(86, 12).Value = 15250  # L86: was 12133.333
$ws.Cells.Item(86, 14).Value = -17496  # N86: was -14379.333

# Sheet ALC, row 89 (Leve Item ID 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 11399.8  # H89: was 10483.167
$ws.Cells.Item(89, 10).Value = 15250  # J89: was 12133.333
$ws.Cells.Item(89, 12).Value = 76250  # L89: was 60666.665
$ws.Cells.Item(89, 14).Value = -87482  # N89: was -71898.66500000001

# Sheet ALC, row 116 (Leve Item ID 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 6215.364  # H116: was 6065.4614
$ws.Cells.Item(116, 9).Value = 6423  # I116: was 6264.8887
$ws.Cells.Item(116, 10).Value = 5661.6665  # J116: was 5616.75
$ws.Cells.Item(116, 11).Value = 6423  # K116: was 6264.8887
$ws.Cells.Item(116, 12).Value = 5661.6665  # L116: was 5616.75
$ws.Cells.Item(116, 13).Value = -2981  # M116: was -2822.8887
$ws.Cells.Item(116, 14).Value = -12545.6665  # N116: was -12500.75

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 19574.75  # H132: was 25675.5
$ws.Cells.Item(132, 9).Value = 17032.584  # I132: was 24900.25
$ws.Cells.Item(132, 10).Value = 27201.25  # J132: was 27226
$ws.Cells.Item(132, 11).Value = 51097.75199999999  # K132: was 74700.75
$ws.Cells.Item(132, 12).Value = 81603.75  # L132: was 81678
$ws.Cells.Item(132, 13).Value = -48567.75199999999  # M132: was -72170.75
$ws.Cells.Item(132, 14).Value = -86663.75  # N132: was -86738

# Sheet ALC, row 141 (Leve Item ID 44161)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 7061.8335  # H141: was 7444.727
$ws.Cells.Item(141, 9).Value = 6474.3  # I141: was 6877
$ws.Cells.Item(141, 11).Value = 19422.9  # K141: was 20631
$ws.Cells.Item(141, 13).Value = -14242.9  # M141: was -15451

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2549.948  # H32: was 2572.2947
$ws.Cells.Item(32, 9).Value = 2041.2675  # I32: was 2047.9186
$ws.Cells.Item(32, 10).Value = 6924.6  # J32: was 7583
$ws.Cells.Item(32, 11).Value = 2041.2675  # K32: was 2047.9186
$ws.Cells.Item(32, 12).Value = 6924.6  # L32: was 7583
$ws.Cells.Item(32, 13).Value = -1754.2675  # M32: was -1760.9186
$ws.Cells.Item(32, 14).Value = -7498.6  # N32: was -8157

# Sheet ARM, row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 85525.67999999999  # H61: was 88519.3
$ws.Cells.Item(61, 10).Value = 181479.92  # J61: was 196211.75
$ws.Cells.Item(61, 12).Value = 181479.92  # L61: was 196211.75
$ws.Cells.Item(61, 14).Value = -181903.92  # N61: was -196635.75

# Sheet ARM, row 74 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 21217.904  # H74: was 22187.7
$ws.Cells.Item(74, 9).Value = 2288.3333  # I74: was 2330.7273
$ws.Cells.Item(74, 11).Value = 2288.3333  # K74: was 2330.7273
$ws.Cells.Item(74, 13).Value = -1414.3333  # M74: was -1456.7273

# Sheet ARM, row 77 (Leve Item ID 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 21217.904  # H77: was 22187.7
$ws.Cells.Item(77, 9).Value = 2288.3333  # I77: was 2330.7273
$ws.Cells.Item(77, 11).Value = 11441.6665  # K77: was 11653.6365
$ws.Cells.Item(77, 13).Value = -7073.666499999999  # M77: was -7285.636500000001

# Sheet ARM, row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 4644.7334  # H110: was 5249.6924
$ws.Cells.Item(110, 9).Value = 6145.9  # I110: was 7504.25
$ws.Cells.Item(110, 11).Value = 6145.9  # K110: was 7504.25
$ws.Cells.Item(110, 13).Value = -4100.9  # M110: was -5459.25

# Sheet ARM, row 111 (Leve Item ID 25813)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(111, 8).Value = 33200  # H111: was 39333.332
$ws.Cells.Item(111, 10).Value = 33200  # J111: was 39333.332
$ws.Cells.Item(111, 12).Value = 33200  # L111: was 39333.332
$ws.Cells.Item(111, 14).Value = -41380  # N111: was -47513.332

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 13201.881  # H132: was 14187.923
$ws.Cells.Item(132, 9).Value = 3093.2693  # I132: was 3446.739
$ws.Cells.Item(132, 11).Value = 9279.8079  # K132: was 10340.217
$ws.Cells.Item(132, 13).Value = -6749.8079  # M132: was -7810.217000000001

# Sheet ARM, row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 85525.67999999999  # H136: was 88519.3
$ws.Cells.Item(136, 10).Value = 181479.92  # J136: was 196211.75
$ws.Cells.Item(136, 12).Value = 544439.76  # L136: was 588635.25
$ws.Cells.Item(136, 14).Value = -549539.76  # N136: was -593735.25

# Sheet BSM, row 134 (Leve Item ID 43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 33617.08  # H134: was 34479.92
$ws.Cells.Item(134, 9).Value = 34689.934  # I134: was 35818.633
$ws.Cells.Item(134, 11).Value = 104069.802  # K134: was 107455.899
$ws.Cells.Item(134, 13).Value = -101534.802  # M134: was -104920.899

# Sheet BSM, row 138 (Leve Item ID 42308)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(138, 8).Value = 78123.125  # H138: was 77221.44500000001
$ws.Cells.Item(138, 10).Value = 78123.125  # J138: was 77221.44500000001
$ws.Cells.Item(138, 12).Value = 78123.125  # L138: was 77221.44500000001
$ws.Cells.Item(138, 14).Value = -88403.125  # N138: was -87501.44500000001

# Sheet BSM, row 140 (Leve Item ID 42471)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 79372.125  # H140: was 78397.89999999999
$ws.Cells.Item(140, 10).Value = 79372.125  # J140: was 78397.89999999999
$ws.Cells.Item(140, 12).Value = 79372.125  # L140: was 78397.89999999999
$ws.Cells.Item(140, 14).Value = -89732.125  # N140: was -88757.89999999999

# Sheet CRP, row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 37229.95  # H31: was 39084.156
$ws.Cells.Item(31, 10).Value = 38631.453  # J31: was 42294.6
$ws.Cells.Item(31, 12).Value = 38631.453  # L31: was 42294.6
$ws.Cells.Item(31, 14).Value = -39221.453  # N31: was -42884.6

# Sheet CRP, row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 37229.95  # H34: was 39084.156
$ws.Cells.Item(34, 10).Value = 38631.453  # J34: was 42294.6
$ws.Cells.Item(34, 12).Value = 38631.453  # L34: was 42294.6
$ws.Cells.Item(34, 14).Value = -39035.453  # N34: was -42698.6

# Sheet CRP, row 38 (Leve Item ID 1637)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(38, 8).Value = 3800  # H38: was 5999.3335
$ws.Cells.Item(38, 9).Value = 3800  # I38: was 5999.3335
$ws.Cells.Item(38, 11).Value = 3800  # K38: was 5999.3335
$ws.Cells.Item(38, 13).Value = -3423  # M38: was -5622.3335

# Sheet CRP, row 46 (Leve Item ID 1637)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(46, 8).Value = 3800  # H46: was 5999.3335
$ws.Cells.Item(46, 9).Value = 3800  # I46: was 5999.3335
$ws.Cells.Item(46, 11).Value = 3800  # K46: was 5999.3335
$ws.Cells.Item(46, 13).Value = -3589  # M46: was -5788.3335

# Sheet CRP, row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 14617.5  # H58: was 15191.046
$ws.Cells.Item(58, 9).Value = 5545.794  # I58: was 5767.4375
$ws.Cells.Item(58, 11).Value = 5545.794  # K58: was 5767.4375
$ws.Cells.Item(58, 13).Value = -5342.794  # M58: was -5564.4375

# Sheet CRP, row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 8302.333000000001  # H132: was 9562.5
$ws.Cells.Item(132, 9).Value = 2987.8572  # I132: was 3301.6667
$ws.Cells.Item(132, 10).Value = 18931.285  # J132: was 22084.166
$ws.Cells.Item(132, 11).Value = 8963.571599999999  # K132: was 9905.000100000001
$ws.Cells.Item(132, 12).Value = 56793.855  # L132: was 66252.49800000001
$ws.Cells.Item(132, 13).Value = -6433.571599999999  # M132: was -7375.000100000001
$ws.Cells.Item(132, 14).Value = -61853.855  # N132: was -71312.49800000001

# Sheet CRP, row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 7779.1665  # H134: was 9531.053
$ws.Cells.Item(134, 9).Value = 2163.0588  # I134: was 2596.8333
$ws.Cells.Item(134, 11).Value = 6489.176399999999  # K134: was 7790.499899999999
$ws.Cells.Item(134, 13).Value = -3954.176399999999  # M134: was -5255.499899999999

# Sheet CRP, row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 14617.5  # H136: was 15191.046
$ws.Cells.Item(136, 9).Value = 5545.794  # I136: was 5767.4375
$ws.Cells.Item(136, 11).Value = 16637.382  # K136: was 17302.3125
$ws.Cells.Item(136, 13).Value = -14087.382  # M136: was -14752.3125

# Sheet CUL, row 2 (Leve Item ID 4847)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 275.78946  # H2: was 275.8421
$ws.Cells.Item(2, 10).Value = 151.375  # J2: was 151.5
$ws.Cells.Item(2, 12).Value = 908.25  # L2: was 909
$ws.Cells.Item(2, 14).Value = -1134.25  # N2: was -1135

# Sheet CUL, row 80 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 12410.167  # H80: was 13057.8
$ws.Cells.Item(80, 9).Value = 1397  # I80: was 1500
$ws.Cells.Item(80, 10).Value = 14612.8  # J80: was 15947.25
$ws.Cells.Item(80, 11).Value = 4191  # K80: was 4500
$ws.Cells.Item(80, 12).Value = 43838.39999999999  # L80: was 47841.75
$ws.Cells.Item(80, 13).Value = -3255  # M80: was -3564
$ws.Cells.Item(80, 14).Value = -45710.39999999999  # N80: was -49713.75

# Sheet CUL, row 83 (Leve Item ID 12890)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 12410.167  # H83: was 13057.8
$ws.Cells.Item(83, 9).Value = 1397  # I83: was 1500
$ws.Cells.Item(83, 10).Value = 14612.8  # J83: was 15947.25
$ws.Cells.Item(83, 11).Value = 12573  # K83: was 13500
$ws.Cells.Item(83, 12).Value = 131515.2  # L83: was 143525.25
$ws.Cells.Item(83, 13).Value = -7893  # M83: was -8820
$ws.Cells.Item(83, 14).Value = -140875.2  # N83: was -152885.25

# Sheet CUL, row 113 (Leve Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1041.88  # H113: was 984.89655
$ws.Cells.Item(113, 9).Value = 1087.5  # I113: was 919.2308
$ws.Cells.Item(113, 10).Value = 1020.41174  # J113: was 1038.25
$ws.Cells.Item(113, 11).Value = 3262.5  # K113: was 2757.6924
$ws.Cells.Item(113, 12).Value = 3061.23522  # L113: was 3114.75
$ws.Cells.Item(113, 13).Value = -1092.5  # M113: was -587.6923999999999
$ws.Cells.Item(113, 14).Value = -7401.23522  # N113: was -7454.75

# Sheet CUL, row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1464.27  # H131: was 1470.73
$ws.Cells.Item(131, 9).Value = 1035.5714  # I131: was 1045.8334
$ws.Cells.Item(131, 10).Value = 1496.5376  # J131: was 1497.8511
$ws.Cells.Item(131, 11).Value = 3106.7142  # K131: was 3137.5002
$ws.Cells.Item(131, 12).Value = 4489.612800000001  # L131: was 4493.5533
$ws.Cells.Item(131, 13).Value = 1933.2858  # M131: was 1902.4998
$ws.Cells.Item(131, 14).Value = -14569.6128  # N131: was -14573.5533

# Sheet CUL, row 134 (Leve Item ID 44074)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value = 6475.7256  # H134: was 6593.38
$ws.Cells.Item(134, 9).Value = 2974.5  # I134: was 3466.2222
$ws.Cells.Item(134, 10).Value = 7127.116  # J134: was 7279.829
$ws.Cells.Item(134, 11).Value = 8923.5  # K134: was 10398.6666
$ws.Cells.Item(134, 12).Value = 21381.348  # L134: was 21839.487
$ws.Cells.Item(134, 13).Value = -3853.5  # M134: was -5328.6666
$ws.Cells.Item(134, 14).Value = -31521.348  # N134: was -31979.487

# Sheet CUL, row 139 (Leve Item ID 44102)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 17987.555  # H139: was 15398.818
$ws.Cells.Item(139, 9).Value = 19486  # I139: was 16338.7
$ws.Cells.Item(139, 11).Value = 58458  # K139: was 49016.10000000001
$ws.Cells.Item(139, 13).Value = -53318  # M139: was -43876.10000000001

# Sheet GSM, row 43 (Leve Item ID 4218)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2975.7083  # H43: was 3016.68
$ws.Cells.Item(43, 9).Value = 2563.5625  # I43: was 2576.0625
$ws.Cells.Item(43, 11).Value = 2563.5625  # K43: was 2576.0625
$ws.Cells.Item(43, 13).Value = -2412.5625  # M43: was -2425.0625

# Sheet GSM, row 107 (Leve Item ID 27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1659.8  # H107: was 1899.75
$ws.Cells.Item(107, 9).Value = 1714.1428  # I107: was 1966.5
$ws.Cells.Item(107, 10).Value = 1533  # J107: was 1699.5
$ws.Cells.Item(107, 11).Value = 1714.1428  # K107: was 1966.5
$ws.Cells.Item(107, 12).Value = 1533  # L107: was 1699.5
$ws.Cells.Item(107, 13).Value = 205.8571999999999  # M107: was -46.5
$ws.Cells.Item(107, 14).Value = -5373  # N107: was -5539.5

# Sheet GSM, row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 10324.158  # H126: was 11177.117
$ws.Cells.Item(126, 9).Value = 13712.223  # I126: was 15213.875
$ws.Cells.Item(126, 10).Value = 7274.9  # J126: was 7588.8887
$ws.Cells.Item(126, 11).Value = 41136.669  # K126: was 45641.625
$ws.Cells.Item(126, 12).Value = 21824.7  # L126: was 22766.6661
$ws.Cells.Item(126, 13).Value = -38666.669  # M126: was -43171.625
$ws.Cells.Item(126, 14).Value = -26764.7  # N126: was -27706.6661

# Sheet LTW, row 70 (Leve Item ID 10811)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(70, 8).Value = 0  # H70: was 13672.5
$ws.Cells.Item(70, 10).Value = 0  # J70: was 13672.5
$ws.Cells.Item(70, 12).Value = 0  # L70: was 13672.5
$ws.Cells.Item(70, 14).ClearContents()  # N70: was -14212.5

# Sheet LTW, row 73 (Leve Item ID 10811)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(73, 8).Value = 0  # H73: was 13672.5
$ws.Cells.Item(73, 10).Value = 0  # J73: was 13672.5
$ws.Cells.Item(73, 12).Value = 0  # L73: was 13672.5
$ws.Cells.Item(73, 14).ClearContents()  # N73: was -15544.5

# Sheet LTW, row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 14670.75  # H132: was 17752.076
$ws.Cells.Item(132, 9).Value = 2099.2222  # I132: was 2489.6667
$ws.Cells.Item(132, 11).Value = 6297.6666  # K132: was 7469.000100000001
$ws.Cells.Item(132, 13).Value = -3767.6666  # M132: was -4939.000100000001

# Sheet LTW, row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 76182.164  # H136: was 76082
$ws.Cells.Item(136, 10).Value = 15937.5  # J136: was 15749.6875
$ws.Cells.Item(136, 12).Value = 47812.5  # L136: was 47249.0625
$ws.Cells.Item(136, 14).Value = -52912.5  # N136: was -52349.0625

# Sheet WVR, row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4735.913  # H132: was 4916
$ws.Cells.Item(132, 9).Value = 1509.973  # I132: was 1552.0286
$ws.Cells.Item(132, 11).Value = 4529.919  # K132: was 4656.085800000001
$ws.Cells.Item(132, 13).Value = -1999.919  # M132: was -2126.085800000001

# Sheet WVR, row 136 (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 14856.2  # H136: was 15398.76
$ws.Cells.Item(136, 9).Value = 1574.7333  # I136: was 1644.0714
$ws.Cells.Item(136, 10).Value = 34778.4  # J136: was 32904.727
$ws.Cells.Item(136, 11).Value = 4724.199900000001  # K136: was 4932.2142
$ws.Cells.Item(136, 12).Value = 104335.2  # L136: was 98714.181
$ws.Cells.Item(136, 13).Value = -2174.199900000001  # M136: was -2382.2142
$ws.Cells.Item(136, 14).Value = -109435.2  # N136: was -103814.181
